# FormattedCarparkRates.xlsx refactor:
#   "all attributes are now camelCase" -> rename the Day value
#   "sunday_public_holiday" to "sundayPublicHoliday" everywhere it is used.
#
# Renaming the shared-string's text (rather than deleting/re-adding it)
# naturally makes Excel drop the old shared-string-table slot and append
# the freshly-named string at the end of the table on save, which is
# exactly the shared-string reshuffle seen in the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDay = "sunday_public_holiday"
$newDay = "sundayPublicHoliday"

$usedRows = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $usedRows; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq $oldDay) {
        $cell.Value = $newDay
    }
}

# Restore the author's last on-screen selection/scroll position.
$ws.Range("B18").Select()
